$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of expense data
$ws.Range("A3").Value = 43178
$ws.Range("B3").Value = "HM-10"
$ws.Range("C3").Value = 115000

$ws.Range("A4").Value = 43178
$ws.Range("B4").Value = "Pin LIPO 500mAh"
$ws.Range("C4").Value = 55000

$ws.Range("A5").Value = 43179
$ws.Range("B5").Value = "Mạch sạc pin LIPO"
$ws.Range("C5").Value = 15000

$ws.Range("A6").Value = 43181
$ws.Range("B6").Value = "HM-10"
$ws.Range("C6").Value = 115000

$ws.Range("A3").NumberFormat = "mm-dd-yy"
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("A3").Copy()
$ws.Range("A4:A6").PasteSpecial(-4122)

$ws.Range("B4").Select()
